$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the style of A20:B20 into A21:B21 so the new row matches the
# existing "Project"/"File" columns formatting (style index 4).
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A21").Value = "Disco.Localization.Resources"
$ws.Range("B21").Value = "Strings"
$ws.Range("C21").Value = "ArgumentOutOfRangeException_invalid_network_port"
$ws.Range("E21").Value = "Invalid network port number ""{0}"". The port number must be an integer value between 0 ~ 65535."
$ws.Range("G21").Value = "Invalid network port number ""{0}"". The port number must be an integer value between 0 ~ 65535."
$ws.Range("I21").Value = "无效的网络端口号“{0}”。端口号必须是一个介于 0 ~ 65535 之间的整型值。"

# Move the active selection to the new last cell, mirroring how Excel
# leaves the cursor on the last-entered row.
$ws.Range("I21").Select()
